# Auto-generated Excel COM-interop script
# Applies cell value updates to the Adamantoise_Profits workbook (sheets ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR)
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value = 1000.54095
$ws.Range("I15").Value = 1000.54095
$ws.Range("K15").Value = 3001.62285
$ws.Range("M15").Value = -2832.62285
$ws.Range("H16").Value = 0
$ws.Range("I16").Value = 0
$ws.Range("K16").Value = 0
$ws.Range("M16").ClearContents()
$ws.Range("H51").Value = 3222.44
$ws.Range("I51").Value = 4642.5835
$ws.Range("K51").Value = 4642.5835
$ws.Range("M51").Value = -4158.5835
$ws.Range("H74").Value = 7195.143
$ws.Range("I74").Value = 5303.5557
$ws.Range("K74").Value = 5303.5557
$ws.Range("M74").Value = -4367.5557
$ws.Range("H77").Value = 7195.143
$ws.Range("I77").Value = 5303.5557
$ws.Range("K77").Value = 26517.7785
$ws.Range("M77").Value = -21837.7785
$ws.Range("H100").Value = 1942.52
$ws.Range("I100").Value = 1754.7858
$ws.Range("J100").Value = 2181.4546
$ws.Range("K100").Value = 1754.7858
$ws.Range("L100").Value = 2181.4546
$ws.Range("M100").Value = -1213.7858
$ws.Range("N100").Value = -3263.4546
$ws.Range("H106").Value = 4450971.5
$ws.Range("I106").Value = 4763184
$ws.Range("K106").Value = 4763184
$ws.Range("M106").Value = -4762553
$ws.Range("H112").Value = 113787.664
$ws.Range("I112").Value = 1187.5
$ws.Range("J112").Value = 145959.14
$ws.Range("K112").Value = 3562.5
$ws.Range("L112").Value = 437877.42
$ws.Range("M112").Value = -2454.5
$ws.Range("N112").Value = -440093.42
$ws.Range("H132").Value = 1834.4231
$ws.Range("I132").Value = 1927.5834
$ws.Range("J132").Value = 716.5
$ws.Range("K132").Value = 5782.7502
$ws.Range("L132").Value = 2149.5
$ws.Range("M132").Value = -3252.7502
$ws.Range("N132").Value = -7209.5
$ws.Range("H133").Value = 0
$ws.Range("J133").Value = 0
$ws.Range("L133").Value = 0
$ws.Range("N133").ClearContents()
$ws.Range("H138").Value = 3644.57
$ws.Range("I138").Value = 3239.4
$ws.Range("J138").Value = 3665.8948
$ws.Range("K138").Value = 9718.200000000001
$ws.Range("L138").Value = 10997.6844
$ws.Range("M138").Value = -4578.200000000001
$ws.Range("N138").Value = -21277.6844

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 17615132
$ws.Range("I32").Value = 20313598
$ws.Range("K32").Value = 20313598
$ws.Range("M32").Value = -20313311
$ws.Range("H37").Value = 35116.57
$ws.Range("J37").Value = 99999
$ws.Range("L37").Value = 99999
$ws.Range("N37").Value = -100545
$ws.Range("H45").Value = 3720.8928
$ws.Range("I45").Value = 3409.3
$ws.Range("J45").Value = 4499.875
$ws.Range("K45").Value = 3409.3
$ws.Range("L45").Value = 4499.875
$ws.Range("M45").Value = -3032.3
$ws.Range("N45").Value = -5253.875
$ws.Range("H74").Value = 3649.5715
$ws.Range("I74").Value = 3649.5715
$ws.Range("K74").Value = 3649.5715
$ws.Range("M74").Value = -2775.5715
$ws.Range("H77").Value = 3649.5715
$ws.Range("I77").Value = 3649.5715
$ws.Range("K77").Value = 18247.8575
$ws.Range("M77").Value = -13879.8575
$ws.Range("H122").Value = 2856.6365
$ws.Range("J122").Value = 3958.6667
$ws.Range("L122").Value = 11876.0001
$ws.Range("N122").Value = -16776.0001

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 3076.1333
$ws.Range("I99").Value = 3085.7273
$ws.Range("K99").Value = 3085.7273
$ws.Range("M99").Value = -1587.7273
$ws.Range("H134").Value = 2153659.2
$ws.Range("I134").Value = 2471731
$ws.Range("K134").Value = 7415193
$ws.Range("M134").Value = -7412658

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 6189.081
$ws.Range("I22").Value = 3128.9285
$ws.Range("J22").Value = 15709.556
$ws.Range("K22").Value = 3128.9285
$ws.Range("L22").Value = 15709.556
$ws.Range("M22").Value = -2778.9285
$ws.Range("N22").Value = -16409.556
$ws.Range("H31").Value = 1914.8113
$ws.Range("I31").Value = 1336.5652
$ws.Range("J31").Value = 2358.1333
$ws.Range("K31").Value = 1336.5652
$ws.Range("L31").Value = 2358.1333
$ws.Range("M31").Value = -1041.5652
$ws.Range("N31").Value = -2948.1333
$ws.Range("H34").Value = 1914.8113
$ws.Range("I34").Value = 1336.5652
$ws.Range("J34").Value = 2358.1333
$ws.Range("K34").Value = 1336.5652
$ws.Range("L34").Value = 2358.1333
$ws.Range("M34").Value = -1134.5652
$ws.Range("N34").Value = -2762.1333
$ws.Range("H58").Value = 2462
$ws.Range("I58").Value = 812.58826
$ws.Range("K58").Value = 812.58826
$ws.Range("M58").Value = -609.58826
$ws.Range("H97").Value = 110000
$ws.Range("J97").Value = 110000
$ws.Range("L97").Value = 110000
$ws.Range("N97").Value = -111982
$ws.Range("H122").Value = 4612.5
$ws.Range("I122").Value = 2607.5625
$ws.Range("K122").Value = 7822.6875
$ws.Range("M122").Value = -5372.6875
$ws.Range("H136").Value = 2462
$ws.Range("I136").Value = 812.58826
$ws.Range("K136").Value = 2437.76478
$ws.Range("M136").Value = 112.23522

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 58.4
$ws.Range("I2").Value = 67.5
$ws.Range("K2").Value = 405
$ws.Range("M2").Value = -292
$ws.Range("H92").Value = 2748.5
$ws.Range("I92").Value = 0
$ws.Range("J92").Value = 2748.5
$ws.Range("K92").Value = 0
$ws.Range("L92").Value = 8245.5
$ws.Range("M92").ClearContents()
$ws.Range("N92").Value = -10741.5
$ws.Range("H107").Value = 1282.9524
$ws.Range("I107").Value = 1015.9167
$ws.Range("K107").Value = 3047.7501
$ws.Range("M107").Value = -1127.7501
$ws.Range("H121").Value = 9495.666999999999
$ws.Range("I121").Value = 417.5
$ws.Range("J121").Value = 18573.834
$ws.Range("K121").Value = 1252.5
$ws.Range("L121").Value = 55721.50199999999
$ws.Range("M121").Value = 57.5
$ws.Range("N121").Value = -58341.50199999999

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H3").Value = 4667666.5
$ws.Range("J3").Value = 0
$ws.Range("L3").Value = 0
$ws.Range("N3").ClearContents()
$ws.Range("H10").Value = 1661.3334
$ws.Range("J10").Value = 1661.3334
$ws.Range("L10").Value = 1661.3334
$ws.Range("N10").Value = -1999.3334
$ws.Range("H12").Value = 2750
$ws.Range("I12").Value = 3000
$ws.Range("J12").Value = 2500
$ws.Range("K12").Value = 3000
$ws.Range("L12").Value = 2500
$ws.Range("M12").Value = -2860
$ws.Range("N12").Value = -2780
$ws.Range("H102").Value = 2342.3845
$ws.Range("I102").Value = 2168.7273
$ws.Range("K102").Value = 2168.7273
$ws.Range("M102").Value = -546.7273
$ws.Range("H122").Value = 1443.4688
$ws.Range("I122").Value = 1256.4762
$ws.Range("K122").Value = 3769.4286
$ws.Range("M122").Value = -1319.4286
$ws.Range("H126").Value = 3571.5264
$ws.Range("I126").Value = 3087.625
$ws.Range("J126").Value = 3923.4546
$ws.Range("K126").Value = 9262.875
$ws.Range("L126").Value = 11770.3638
$ws.Range("M126").Value = -6792.875
$ws.Range("N126").Value = -16710.3638
$ws.Range("H132").Value = 2937.182
$ws.Range("I132").Value = 2333
$ws.Range("K132").Value = 6999
$ws.Range("M132").Value = -4469

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 2404
$ws.Range("I40").Value = 2404
$ws.Range("K40").Value = 2404
$ws.Range("M40").Value = -2268
$ws.Range("H61").Value = 6177.5557
$ws.Range("I61").Value = 1919.8
$ws.Range("J61").Value = 11499.75
$ws.Range("K61").Value = 1919.8
$ws.Range("L61").Value = 11499.75
$ws.Range("M61").Value = -1717.8
$ws.Range("N61").Value = -11903.75
$ws.Range("H100").Value = 2833.3333
$ws.Range("J100").Value = 2833.3333
$ws.Range("L100").Value = 2833.3333
$ws.Range("N100").Value = -3915.3333
$ws.Range("H113").Value = 6177.5557
$ws.Range("I113").Value = 1919.8
$ws.Range("J113").Value = 11499.75
$ws.Range("K113").Value = 1919.8
$ws.Range("L113").Value = 11499.75
$ws.Range("M113").Value = 250.2
$ws.Range("N113").Value = -15839.75
$ws.Range("H131").Value = 140999
$ws.Range("J131").Value = 140999
$ws.Range("L131").Value = 140999
$ws.Range("N131").Value = -151079
$ws.Range("H132").Value = 7414.5674
$ws.Range("I132").Value = 7024.3794
$ws.Range("K132").Value = 21073.1382
$ws.Range("M132").Value = -18543.1382
$ws.Range("H136").Value = 51996.43
$ws.Range("I136").Value = 59830
$ws.Range("J136").Value = 4995
$ws.Range("K136").Value = 179490
$ws.Range("L136").Value = 14985
$ws.Range("M136").Value = -176940
$ws.Range("N136").Value = -20085

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H45").Value = 13612.75
$ws.Range("J45").Value = 13612.75
$ws.Range("L45").Value = 13612.75
$ws.Range("N45").Value = -14594.75
$ws.Range("H122").Value = 5714.1665
$ws.Range("H132").Value = 2854.0645
$ws.Range("I132").Value = 2694.9583
$ws.Range("K132").Value = 8084.874899999999
$ws.Range("M132").Value = -5554.874899999999
